$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained 3 rows of earlier historical data (2019-11-26 .. 2019-11-28),
# inserted right after the existing first data row, pushing everything else
# down by 3 rows (old row 2 -> new row 5, ..., old row 74 -> new row 77).
#
# Rows are inserted below row 2 (a plain, unstyled data row) rather than
# directly below the bold header row, so the newly created rows don't
# inherit the header's bold/centered style. Row 2's original content is
# then copied down to row 5, and rows 2-4 are overwritten with the new
# historical data.
$ws.Rows.Item(3).Resize(3).Insert()

$ws.Range("A2:I2").Copy($ws.Range("A5:I5"))

function Set-TextValue($range, $text) {
    # Force literal text storage (avoids Excel auto-coercing numeric-looking
    # strings like "0215" or date-looking strings like "2019-11-26" into
    # numbers/dates), then drop the temporary Text format so the cell keeps
    # the workbook's default (unstyled) look.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("A2").Value = 1574726400
Set-TextValue $ws.Range("B2") "2019-11-26"
Set-TextValue $ws.Range("C2") "0215"
Set-TextValue $ws.Range("D2") "SLVEST"
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.76
$ws.Range("G2").Value = 0.595
$ws.Range("H2").Value = 0.755
$ws.Range("I2").Value = 262939400

$ws.Range("A3").Value = 1574812800
Set-TextValue $ws.Range("B3") "2019-11-27"
Set-TextValue $ws.Range("C3") "0215"
Set-TextValue $ws.Range("D3") "SLVEST"
$ws.Range("E3").Value = 0.77
$ws.Range("F3").Value = 0.85
$ws.Range("G3").Value = 0.76
$ws.Range("H3").Value = 0.82
$ws.Range("I3").Value = 134962300

$ws.Range("A4").Value = 1574899200
Set-TextValue $ws.Range("B4") "2019-11-28"
Set-TextValue $ws.Range("C4") "0215"
Set-TextValue $ws.Range("D4") "SLVEST"
$ws.Range("E4").Value = 0.83
$ws.Range("F4").Value = 0.83
$ws.Range("G4").Value = 0.74
$ws.Range("H4").Value = 0.755
$ws.Range("I4").Value = 58524100
